# Se sube reporte html
#
# The "Expected Content" column is split into six columns
# (Expected Content_1 .. Expected Content_6) so each test row can record
# several expected JSON fragments; a couple of rows are also toggled on/off
# and refreshed with the real API responses captured while producing the
# new HTML report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row: "Expected Content" -> "Expected Content_1", plus five new
#    headers for the columns that follow it.
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = "Expected Content_1"
$ws.Range("H1").Value = "Expected Content_2"
$ws.Range("I1").Value = "Expected Content_3"
$ws.Range("J1").Value = "Expected Content_4"
$ws.Range("K1").Value = "Expected Content_5"
$ws.Range("L1").Value = "Expected Content_6"

# ---------------------------------------------------------------------------
# 2. Value edits on the existing data rows (2-6).
# ---------------------------------------------------------------------------
# Row 2 - login with valid credentials: no longer checks Expected Content_1.
$ws.Range("G2").Value = ""

# Row 3 - login with missing password: now executed (NO -> SI) and records
# the API's error payload in Expected Content_2.
$ws.Range("A3").Value = "SI"
$ws.Range("H3").Value = "{
    ""error"": ""Missing password""
}"

# Row 4 - create user: now executed (NO -> SI).
$ws.Range("A4").Value = "SI"

# Row 5 - GET todo: Expected Content_1 becomes the real pretty-printed
# response body, and the id it carries is captured separately.
$ws.Range("G5").Value = "{
    ""id"": 1
}"
$ws.Range("H5").Value = "id"
$ws.Range("I5").Value = 1
$ws.Rows.Item(5).RowHeight = 43.5

# Row 6 - create post: no longer executed (SI -> NO) and Expected Content_1
# is refreshed with the actual response.
$ws.Range("A6").Value = "NO"
$ws.Range("G6").Value = "{""id"":1}"

# ---------------------------------------------------------------------------
# 3. Formatting: extend the banded look of column G into the new H:L
#    columns, row by row, using each row's own formatting as the template.
# ---------------------------------------------------------------------------
$ws.Range("G1").Copy()
$ws.Range("H1:L1").PasteSpecial(-4122)

$ws.Range("G2").Copy()
$ws.Range("H2:L2").PasteSpecial(-4122)

$ws.Range("G3").Copy()
$ws.Range("H3:L3").PasteSpecial(-4122)

$ws.Range("G4").Copy()
$ws.Range("H4:L4").PasteSpecial(-4122)

# Row 5: H5:L5 use the regular banded style (same as G2/G3/G4); G5 itself
# switches to the plain interior style used by the rest of row 5 (e.g. F5).
$ws.Range("G2").Copy()
$ws.Range("H5:L5").PasteSpecial(-4122)
$ws.Range("F5").Copy()
$ws.Range("G5").PasteSpecial(-4122)

# Row 6: H6:J6 use the regular banded style; K6:L6 keep the plain interior
# style that column G already has on this row (matches the source file).
$ws.Range("G2").Copy()
$ws.Range("H6:J6").PasteSpecial(-4122)
$ws.Range("G6").Copy()
$ws.Range("K6:L6").PasteSpecial(-4122)

# Blank rows 7-46: copy each row's own column-G formatting across H:L.
for ($r = 7; $r -le 46; $r++) {
    $ws.Range("G$r").Copy()
    $ws.Range("H$r`:L$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Column widths: G keeps its visual width but is no longer auto-sized;
#    the new H:L columns pick up that same width.
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 17.17
$ws.Columns.Item(8).ColumnWidth = 17.17
$ws.Columns.Item(9).ColumnWidth = 17.17
$ws.Columns.Item(10).ColumnWidth = 17.17
$ws.Columns.Item(11).ColumnWidth = 17.17
$ws.Columns.Item(12).ColumnWidth = 17.17

# ---------------------------------------------------------------------------
# 5. Cursor ends on I2, matching the author's saved selection.
# ---------------------------------------------------------------------------
$ws.Range("I2").Select()
